$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 677.3125
$ws.Range("J17").Value = 650.587
$ws.Range("L17").Value = 1951.761
$ws.Range("N17").Value = -2287.761

$ws.Range("H33").Value = 238.8
$ws.Range("I33").Value = 201.14285
$ws.Range("K33").Value = 201.14285
$ws.Range("M33").Value = 27.85714999999999

$ws.Range("H43").Value = 4380.769
$ws.Range("I43").Value = 5000
$ws.Range("J43").Value = 4105.5557
$ws.Range("K43").Value = 5000
$ws.Range("L43").Value = 4105.5557
$ws.Range("M43").Value = -4931
$ws.Range("N43").Value = -4243.5557

$ws.Range("H88").Value = 3501.682
$ws.Range("J88").Value = 2493.4119
$ws.Range("L88").Value = 2493.4119
$ws.Range("N88").Value = -3305.4119

$ws.Range("H91").Value = 3501.682
$ws.Range("J91").Value = 2493.4119
$ws.Range("L91").Value = 2493.4119
$ws.Range("N91").Value = -5301.4119

$ws.Range("H129").Value = 289774.28
$ws.Range("I129").Value = 336403.34
$ws.Range("K129").Value = 1009210.02
$ws.Range("M129").Value = -1004210.02

$ws.Range("H131").Value = 10077
$ws.Range("I131").Value = 8287.076999999999
$ws.Range("K131").Value = 24861.231
$ws.Range("M131").Value = -19821.231

$ws.Range("H132").Value = 14045.768
$ws.Range("I132").Value = 1362.6364
$ws.Range("K132").Value = 4087.9092
$ws.Range("M132").Value = -1557.9092

$ws.Range("H138").Value = 3482.8572
$ws.Range("I138").Value = 2126.8333
$ws.Range("J138").Value = 4499.875
$ws.Range("K138").Value = 6380.499899999999
$ws.Range("L138").Value = 13499.625
$ws.Range("M138").Value = -1240.499899999999
$ws.Range("N138").Value = -23779.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13775.069
$ws.Range("I2").Value = 17785
$ws.Range("J2").Value = 4864.1113
$ws.Range("K2").Value = 17785
$ws.Range("L2").Value = 4864.1113
$ws.Range("M2").Value = -17672
$ws.Range("N2").Value = -5090.1113

$ws.Range("H45").Value = 4415.793
$ws.Range("I45").Value = 2301.0625
$ws.Range("J45").Value = 7018.5386
$ws.Range("K45").Value = 2301.0625
$ws.Range("L45").Value = 7018.5386
$ws.Range("M45").Value = -1924.0625
$ws.Range("N45").Value = -7772.5386

$ws.Range("H69").Value = 59998.332
$ws.Range("J69").Value = 59998.332
$ws.Range("L69").Value = 59998.332
$ws.Range("N69").Value = -61496.332

$ws.Range("H72").Value = 59998.332
$ws.Range("J72").Value = 59998.332
$ws.Range("L72").Value = 179994.996
$ws.Range("N72").Value = -187482.996

$ws.Range("H110").Value = 5463.375
$ws.Range("I110").Value = 4682
$ws.Range("K110").Value = 4682
$ws.Range("M110").Value = -2637

$ws.Range("H116").Value = 13775.069
$ws.Range("I116").Value = 17785
$ws.Range("J116").Value = 4864.1113
$ws.Range("K116").Value = 17785
$ws.Range("L116").Value = 4864.1113
$ws.Range("M116").Value = -15491
$ws.Range("N116").Value = -9452.1113

$ws.Range("H122").Value = 4075.6345
$ws.Range("I122").Value = 3343.0322
$ws.Range("J122").Value = 5157.095
$ws.Range("K122").Value = 10029.0966
$ws.Range("L122").Value = 15471.285
$ws.Range("M122").Value = -7579.096600000001
$ws.Range("N122").Value = -20371.285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13775.069
$ws.Range("I3").Value = 17785
$ws.Range("J3").Value = 4864.1113
$ws.Range("K3").Value = 17785
$ws.Range("L3").Value = 4864.1113
$ws.Range("M3").Value = -17671
$ws.Range("N3").Value = -5092.1113

$ws.Range("H24").Value = 397.5
$ws.Range("I24").Value = 397.5
$ws.Range("K24").Value = 397.5
$ws.Range("M24").Value = -162.5

$ws.Range("H99").Value = 20474.334
$ws.Range("I99").Value = 23397.889
$ws.Range("J99").Value = 2933
$ws.Range("K99").Value = 23397.889
$ws.Range("L99").Value = 2933
$ws.Range("M99").Value = -21899.889
$ws.Range("N99").Value = -5929

$ws.Range("H107").Value = 10052.533
$ws.Range("I107").Value = 11478.8
$ws.Range("J107").Value = 7200
$ws.Range("K107").Value = 11478.8
$ws.Range("L107").Value = 7200
$ws.Range("M107").Value = -9558.799999999999
$ws.Range("N107").Value = -11040

$ws.Range("H134").Value = 6450
$ws.Range("I134").Value = 5900
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 17700
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -15165
$ws.Range("N134").Value = -26070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2937.3635
$ws.Range("I16").Value = 2879
$ws.Range("K16").Value = 2879
$ws.Range("M16").Value = -2592

$ws.Range("H25").Value = 3307.5
$ws.Range("I25").Value = 520
$ws.Range("J25").Value = 7953.3335
$ws.Range("K25").Value = 520
$ws.Range("L25").Value = 7953.3335
$ws.Range("M25").Value = -346
$ws.Range("N25").Value = -8301.333500000001

$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H31").Value = 3860.5
$ws.Range("J31").Value = 4874.375
$ws.Range("L31").Value = 4874.375
$ws.Range("N31").Value = -5464.375

$ws.Range("H34").Value = 3860.5
$ws.Range("J34").Value = 4874.375
$ws.Range("L34").Value = 4874.375
$ws.Range("N34").Value = -5278.375

$ws.Range("H113").Value = 2937.3635
$ws.Range("I113").Value = 2879
$ws.Range("K113").Value = 2879
$ws.Range("M113").Value = -709

$ws.Range("H134").Value = 4616.276
$ws.Range("I134").Value = 4640
$ws.Range("K134").Value = 13920
$ws.Range("M134").Value = -11385

$ws.Range("H141").Value = 47500
$ws.Range("J141").Value = 47500
$ws.Range("L141").Value = 47500
$ws.Range("N141").Value = -57860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 2111.25
$ws.Range("I14").Value = 2111.25
$ws.Range("K14").Value = 6333.75
$ws.Range("M14").Value = -6160.75

$ws.Range("H113").Value = 937
$ws.Range("J113").Value = 692.375
$ws.Range("L113").Value = 2077.125
$ws.Range("N113").Value = -6417.125

$ws.Range("H121").Value = 1530.9048
$ws.Range("I121").Value = 547
$ws.Range("J121").Value = 2613.2
$ws.Range("K121").Value = 1641
$ws.Range("L121").Value = 7839.599999999999
$ws.Range("M121").Value = -331
$ws.Range("N121").Value = -10459.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 32499.5
$ws.Range("I52").Value = 30000
$ws.Range("K52").Value = 30000
$ws.Range("M52").Value = -29741

$ws.Range("H70").Value = 106263.63
$ws.Range("I70").Value = 143262.75
$ws.Range("K70").Value = 143262.75
$ws.Range("M70").Value = -142992.75

$ws.Range("H73").Value = 106263.63
$ws.Range("I73").Value = 143262.75
$ws.Range("K73").Value = 143262.75
$ws.Range("M73").Value = -142326.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1479.6666
$ws.Range("I22").Value = 1075.75
$ws.Range("J22").Value = 1626.5454
$ws.Range("K22").Value = 1075.75
$ws.Range("L22").Value = 1626.5454
$ws.Range("M22").Value = -780.75
$ws.Range("N22").Value = -2216.5454

$ws.Range("H27").Value = 1479.6666
$ws.Range("I27").Value = 1075.75
$ws.Range("J27").Value = 1626.5454
$ws.Range("K27").Value = 1075.75
$ws.Range("L27").Value = 1626.5454
$ws.Range("M27").Value = -968.75
$ws.Range("N27").Value = -1840.5454

$ws.Range("H43").Value = 29027.334
$ws.Range("I43").Value = 17137.334
$ws.Range("K43").Value = 17137.334
$ws.Range("M43").Value = -16944.334

$ws.Range("H46").Value = 4256.9644
$ws.Range("J46").Value = 4587.92
$ws.Range("L46").Value = 4587.92
$ws.Range("N46").Value = -4963.92

$ws.Range("H61").Value = 1923.3636
$ws.Range("I61").Value = 1865.7
$ws.Range("K61").Value = 1865.7
$ws.Range("M61").Value = -1663.7

$ws.Range("H113").Value = 1923.3636
$ws.Range("I113").Value = 1865.7
$ws.Range("K113").Value = 1865.7
$ws.Range("M113").Value = 304.3

$ws.Range("H122").Value = 5400
$ws.Range("I122").Value = 3960
$ws.Range("K122").Value = 11880
$ws.Range("M122").Value = -9430

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1018.3125
$ws.Range("I107").Value = 1126.909
$ws.Range("K107").Value = 3380.727
$ws.Range("M107").Value = -1460.727

$ws.Range("H113").Value = 1899.909
$ws.Range("I113").Value = 699.8570999999999
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 2099.5713
$ws.Range("L113").Value = 12000
$ws.Range("M113").Value = 70.42870000000039
$ws.Range("N113").Value = -16340

$ws.Range("H136").Value = 6326.5557
$ws.Range("I136").Value = 6242.375
$ws.Range("K136").Value = 18727.125
$ws.Range("M136").Value = -16177.125
